# Reworking the admin part:
# - Rename the group sheets, dropping the "Группа " ("Group ") prefix
# - Make the first sheet ("А-1-26") the active/selected tab instead of the second one

$wb = $excel.ActiveWorkbook

$wb.Worksheets.Item(1).Name = "А-1-26"
$wb.Worksheets.Item(2).Name = "А-2-26"
$wb.Worksheets.Item(3).Name = "К-1-26"

$wb.Worksheets.Item(1).Activate()
